# Generate Report for Handoff
# Adds two newly-handed-off files (rows 4 & 5) to the Overview sheet and to
# the per-locale detail sheets (zh-cn / de-de), mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$uuid1 = "4b826826-9e3c-453e-9415-f069ad2951bc"
$hash1 = "808f80135784ad36983783656e18e1d384c1cf65"
$uuid2 = "8e04b8e7-4f2d-4ece-a894-cb252dddea8f"
$hash2 = "72702ec6989d6e4f9f7b9a0f8810134f5abb22be"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A4").Value = ($uuid1 + ".md")
$wsOv.Range("B4").Value = "Ready for handoff"
$wsOv.Range("C4").Value = "Ready for handoff"
$wsOv.Range("D4").Value = "2016-48-13 22:48:04"

$wsOv.Range("A5").Value = ($uuid2 + ".md")
$wsOv.Range("B5").Value = "Ready for handoff"
$wsOv.Range("C5").Value = "Ready for handoff"
$wsOv.Range("D5").Value = "2016-48-13 22:48:04"

$wsOv.Hyperlinks.Add($wsOv.Range("A4"), ("https://github.com/OpenLocalizationTest/oltest/blob/ba9b6a6f7b3cf5c0a9a1c47f6b57b1f2f9d8a6f1/e2e/" + $uuid1 + ".md"), "", "", ($uuid1 + ".md"))
$wsOv.Hyperlinks.Add($wsOv.Range("A5"), ("https://github.com/OpenLocalizationTest/oltest/blob/c1d9f4a5e6b7c8d9e0f1a2b3c4d5e6f7a8b9c0d1/e2e/" + $uuid2 + ".md"), "", "", ($uuid2 + ".md"))

$wsOv.Range("A4").Font.Underline = $true
$wsOv.Range("A4").Font.Color = 15570276
$wsOv.Range("A5").Font.Underline = $true
$wsOv.Range("A5").Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = ($uuid1 + ".md")
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = ($uuid1 + "." + $hash1 + ".zh-cn.xlf")
$wsZh.Range("E4").Value = "2016-03-13 22:48:00"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

$wsZh.Range("A5").Value = ($uuid2 + ".md")
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = ($uuid2 + "." + $hash2 + ".zh-cn.xlf")
$wsZh.Range("E5").Value = "2016-03-13 22:48:00"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ("https://github.com/OpenLocalizationTest/oltest/blob/ba9b6a6f7b3cf5c0a9a1c47f6b57b1f2f9d8a6f1/e2e/" + $uuid1 + ".md"), "", "", ($uuid1 + ".md"))
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), ("https://github.com/OpenLocalizationTest/oltest/blob/ba9b6a6f7b3cf5c0a9a1c47f6b57b1f2f9d8a6f1/e2e/" + $uuid1 + ".md"), "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e2f3a4b5c6d7e8f9a0b1c2d3e4f5a6b7c8d9e0f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuid1 + "." + $hash1 + ".zh-cn.xlf"), "", "", ($uuid1 + "." + $hash1 + ".zh-cn.xlf"))

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ("https://github.com/OpenLocalizationTest/oltest/blob/c1d9f4a5e6b7c8d9e0f1a2b3c4d5e6f7a8b9c0d1/e2e/" + $uuid2 + ".md"), "", "", ($uuid2 + ".md"))
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), ("https://github.com/OpenLocalizationTest/oltest/blob/c1d9f4a5e6b7c8d9e0f1a2b3c4d5e6f7a8b9c0d1/e2e/" + $uuid2 + ".md"), "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2f3a4b5c6d7e8f9a0b1c2d3e4f5a6b7c8d9e0f1a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $uuid2 + "." + $hash2 + ".zh-cn.xlf"), "", "", ($uuid2 + "." + $hash2 + ".zh-cn.xlf"))

$wsZh.Range("A4").Font.Underline = $true
$wsZh.Range("A4").Font.Color = 15570276
$wsZh.Range("B4").Font.Underline = $true
$wsZh.Range("B4").Font.Color = 15570276
$wsZh.Range("D4").Font.Underline = $true
$wsZh.Range("D4").Font.Color = 15570276
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("A5").Font.Underline = $true
$wsZh.Range("A5").Font.Color = 15570276
$wsZh.Range("B5").Font.Underline = $true
$wsZh.Range("B5").Font.Color = 15570276
$wsZh.Range("D5").Font.Underline = $true
$wsZh.Range("D5").Font.Color = 15570276
$wsZh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = ($uuid1 + ".md")
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = ($uuid1 + "." + $hash1 + ".de-de.xlf")
$wsDe.Range("E4").Value = "2016-03-13 22:48:04"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

$wsDe.Range("A5").Value = ($uuid2 + ".md")
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = ($uuid2 + "." + $hash2 + ".de-de.xlf")
$wsDe.Range("E5").Value = "2016-03-13 22:48:04"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ("https://github.com/OpenLocalizationTest/oltest/blob/ba9b6a6f7b3cf5c0a9a1c47f6b57b1f2f9d8a6f1/e2e/" + $uuid1 + ".md"), "", "", ($uuid1 + ".md"))
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), ("https://github.com/OpenLocalizationTest/oltest/blob/ba9b6a6f7b3cf5c0a9a1c47f6b57b1f2f9d8a6f1/e2e/" + $uuid1 + ".md"), "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3a4b5c6d7e8f9a0b1c2d3e4f5a6b7c8d9e0f1a2b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuid1 + "." + $hash1 + ".de-de.xlf"), "", "", ($uuid1 + "." + $hash1 + ".de-de.xlf"))

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ("https://github.com/OpenLocalizationTest/oltest/blob/c1d9f4a5e6b7c8d9e0f1a2b3c4d5e6f7a8b9c0d1/e2e/" + $uuid2 + ".md"), "", "", ($uuid2 + ".md"))
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), ("https://github.com/OpenLocalizationTest/oltest/blob/c1d9f4a5e6b7c8d9e0f1a2b3c4d5e6f7a8b9c0d1/e2e/" + $uuid2 + ".md"), "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4b5c6d7e8f9a0b1c2d3e4f5a6b7c8d9e0f1a2b3c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $uuid2 + "." + $hash2 + ".de-de.xlf"), "", "", ($uuid2 + "." + $hash2 + ".de-de.xlf"))

$wsDe.Range("A4").Font.Underline = $true
$wsDe.Range("A4").Font.Color = 15570276
$wsDe.Range("B4").Font.Underline = $true
$wsDe.Range("B4").Font.Color = 15570276
$wsDe.Range("D4").Font.Underline = $true
$wsDe.Range("D4").Font.Color = 15570276
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("A5").Font.Underline = $true
$wsDe.Range("A5").Font.Color = 15570276
$wsDe.Range("B5").Font.Underline = $true
$wsDe.Range("B5").Font.Color = 15570276
$wsDe.Range("D5").Font.Underline = $true
$wsDe.Range("D5").Font.Color = 15570276
$wsDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
